# agregue cuenta 60150 f. alvarez
# The account number 60150 is added to the client list of the
# "RECURRENTE_100K-200K" group, which lives in cell B6 of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B6")
$currentValue = $cell.Value2
$cell.Value = $currentValue + ".60150"

# Reflect the cursor/selection ending on B7, as in the saved workbook.
$ws.Range("B7").Select()
